$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B4 holds the "FilesTab" Neo4j query text (corrected ICDC Breed script).
# The fix drops the "File Type" and "Breed" output columns from the RETURN
# clause - they weren't wanted on the Files tab.
$newFilesQuery = @"
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Staffordshire Bull Terrier'] 
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS ``File Name``,
        coalesce(labels(parent)[0], '') AS ``Association``,
        coalesce(f.file_description, '') AS ``Description``,
        coalesce(f.file_format, '') AS ``Format``,
        coalesce(f.file_size, '') AS ``Size``,
        coalesce(c.case_id, '') AS ``Case ID``,
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS ``Study Code``
"@

$ws.Range("B4").Value = $newFilesQuery

# The wrapped text now takes fewer lines, so the row shrinks to match.
$ws.Rows.Item(4).RowHeight = 217.5

# Update the saved view state: scrolled down with B4 as the active selection.
[void]$ws.Range("B4").Select()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
